# Database_Info.xlsx update:
#  - Fix "Makup" typo -> "Makeup"
#  - Rework the "Department" lookup list in column H (drop the stray
#    "Genre" row, rename "Category/Type" -> "Category/Genre", and add a
#    "Type"/"Type_Lookup" pair at the end of the list)
#  - Document four new tables (Customer, ChartData, Order_Master,
#    Order_Details) in columns I:L, plus their related lookup tables
#  - Extend the "Related Tables" banner row and its merge from A16:H16
#    to A16:L16
#  - Add a new, empty "Sheet2" after "Sheet1" (placeholder for the new
#    data model)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Fix the "Makup" typo in the Makeup table header
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Makeup"

# ---------------------------------------------------------------------
# 2. Column H - "Department" related list: remove the stray "Genre"
#    entry (rows shift up one), rename Category/Type -> Category/Genre,
#    and append Type / Type_Lookup at the bottom in place of the old
#    (incorrect) duplicated Kitchen_Types / Kitchen_Types_Lookup rows.
# ---------------------------------------------------------------------
$ws.Range("H19").Value = "Category/Genre"
$ws.Range("H21").Value = "Author"
$ws.Range("H22").Value = "Author_Lookup"
$ws.Range("H23").Value = "Media"
$ws.Range("H24").Value = "Media_Lookup"
$ws.Range("H25").Value = "Color"
$ws.Range("H26").Value = "Color_Lookup"
$ws.Range("H27").Value = "Size"
$ws.Range("H28").Value = "Size_Lookup"
$ws.Range("H29").Value = "Type"
$ws.Range("H30").Value = "Type_Lookup"
$ws.Range("H31").ClearContents()

# Carry the existing "grey lookup row" look down onto the two new rows
$ws.Range("C20").Copy()
$ws.Range("H29:H30").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. New tables documented in columns I (Customer), J (ChartData),
#    K (Order_Master) and L (Order_Details)
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "Customer"
$ws.Range("I2").Value = "Customer_ID"
$ws.Range("I3").Value = "Customer_FName"
$ws.Range("I4").Value = "Customer_LName"
$ws.Range("I5").Value = "Sex"
$ws.Range("I6").Value = "Birth_Date"
$ws.Range("I7").Value = "Zip_Code"
$ws.Range("I8").Value = "City"
$ws.Range("I9").Value = "State"
$ws.Range("I10").Value = "Street_Number"
$ws.Range("I11").Value = "Street_Name"
$ws.Range("I12").Value = "Marital_Status_Type"

$ws.Range("J1").Value = "ChartData"
$ws.Range("J2").Value = "Data_ID"
$ws.Range("J3").Value = "Year"
$ws.Range("J4").Value = "Order_Count"

$ws.Range("K1").Value = "Order_Master"
$ws.Range("K2").Value = "Order_ID"
$ws.Range("K3").Value = "Order_DateTime"
$ws.Range("K4").Value = "Customer_ID (PK)"

$ws.Range("L1").Value = "Order_Details"
$ws.Range("L2").Value = "Order_Detail_ID"
$ws.Range("L3").Value = "Order_ID"
$ws.Range("L4").Value = "Product_ID"
$ws.Range("L5").Value = "Product_Media_ID"
$ws.Range("L6").Value = "Price"

# Match the bold/centered table-header look already used in H1 for the
# four new table headers
$ws.Range("H1").Copy()
$ws.Range("I1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. "Related Tables" section additions: Customer lookups (I) and a
#    cross reference to the Customer table (K18)
# ---------------------------------------------------------------------
$ws.Range("I18").Value = "Street_Type_Lookup"
$ws.Range("I19").Value = "Degree_Lookup"
$ws.Range("I20").Value = "Income_Lookup"
$ws.Range("I21").Value = "Credit_Card_Type"
$ws.Range("K18").Value = "Customer"

# Highlight the new Street_Type_Lookup entry in yellow
$ws.Range("I18").Interior.Color = 65535

# ---------------------------------------------------------------------
# 5. Extend the "Related Tables" banner row (and its merge) out to L16
# ---------------------------------------------------------------------
$ws.Range("A16:H16").UnMerge()
$ws.Range("H16").Copy()
$ws.Range("I16:L16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A16:L16").Merge()

# ---------------------------------------------------------------------
# 6. Column widths for the new columns (matches the author's layout)
# ---------------------------------------------------------------------
$ws.Range("I1").ColumnWidth = 27.90625
$ws.Range("J1").ColumnWidth = 13.36328125
$ws.Range("K1").ColumnWidth = 21.6328125
$ws.Range("L1").ColumnWidth = 21.7265625

# ---------------------------------------------------------------------
# 7. New blank "Sheet2" placed right after "Sheet1"
# ---------------------------------------------------------------------
$sheet2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$sheet2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 8. Restore the original active sheet / selection
# ---------------------------------------------------------------------
$ws.Select()
$ws.Range("F21").Select()
